$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: prepare formatting for new column I (header) and new rows 51-57 (column A) ---
$ws.Cells.Item(1,8).Copy($ws.Cells.Item(1,9))   # copy H1 style (bold/border/center) to I1
for ($r = 51; $r -le 57; $r++) {
    $ws.Cells.Item(2,1).Copy($ws.Cells.Item($r,1))   # copy A2 style to new A-column label cells
}

# --- Step 2: clear old data cells (values only; keeps styles/labels) so stale entries from the
#     pre-insertion column/row layout do not linger after the shift ---
$ws.Range("B2:I57").ClearContents()

# --- Step 3: rewrite header row (d=1 .. d=10, now including the new d=6 column) ---
$ws.Cells.Item(1,2).Value = "d=1"
$ws.Cells.Item(1,3).Value = "d=2"
$ws.Cells.Item(1,4).Value = "d=3"
$ws.Cells.Item(1,5).Value = "d=4"
$ws.Cells.Item(1,6).Value = "d=5"
$ws.Cells.Item(1,7).Value = "d=6"
$ws.Cells.Item(1,8).Value = "d=7"
$ws.Cells.Item(1,9).Value = "d=10"

# --- Step 4: rewrite column A labels + the single data value for every row (2-57) ---
$ws.Cells.Item(2,1).Value = "ARMA_I(0,1,0)"
$ws.Cells.Item(2,2).Value = 60.71984751580217
$ws.Cells.Item(3,1).Value = "ARMA_I(0,1,1)"
$ws.Cells.Item(3,2).Value = 66.07015436451078
$ws.Cells.Item(4,1).Value = "ARMA_I(0,1,2)"
$ws.Cells.Item(4,2).Value = 81.12193842253824
$ws.Cells.Item(5,1).Value = "ARMA_I(0,10,0)"
$ws.Cells.Item(5,9).Value = 96.18691285694622
$ws.Cells.Item(6,1).Value = "ARMA_I(0,10,1)"
$ws.Cells.Item(6,9).Value = 96.20100698668097
$ws.Cells.Item(7,1).Value = "ARMA_I(0,10,2)"
$ws.Cells.Item(7,9).Value = 96.18082039348931
$ws.Cells.Item(8,1).Value = "ARMA_I(0,2,0)"
$ws.Cells.Item(8,3).Value = 86.09332943401988
$ws.Cells.Item(9,1).Value = "ARMA_I(0,2,1)"
$ws.Cells.Item(9,3).Value = 86.56770807175711
$ws.Cells.Item(10,1).Value = "ARMA_I(0,2,2)"
$ws.Cells.Item(10,3).Value = 86.81184556348657
$ws.Cells.Item(11,1).Value = "ARMA_I(0,3,0)"
$ws.Cells.Item(11,4).Value = 99.21958148675363
$ws.Cells.Item(12,1).Value = "ARMA_I(0,3,1)"
$ws.Cells.Item(12,4).Value = 99.27802673591547
$ws.Cells.Item(13,1).Value = "ARMA_I(0,3,2)"
$ws.Cells.Item(13,4).Value = 99.39171778791176
$ws.Cells.Item(14,1).Value = "ARMA_I(0,4,0)"
$ws.Cells.Item(14,5).Value = 98.93658064173701
$ws.Cells.Item(15,1).Value = "ARMA_I(0,4,1)"
$ws.Cells.Item(15,5).Value = 98.99412345848863
$ws.Cells.Item(16,1).Value = "ARMA_I(0,4,2)"
$ws.Cells.Item(16,5).Value = 98.9458616453217
$ws.Cells.Item(17,1).Value = "ARMA_I(0,5,0)"
$ws.Cells.Item(17,6).Value = 98.47513980388757
$ws.Cells.Item(18,1).Value = "ARMA_I(0,5,1)"
$ws.Cells.Item(18,6).Value = 98.55337953282209
$ws.Cells.Item(19,1).Value = "ARMA_I(0,5,2)"
$ws.Cells.Item(19,6).Value = 98.45312589083835
$ws.Cells.Item(20,1).Value = "ARMA_I(0,6,0)"
$ws.Cells.Item(20,7).Value = 98.11702806779181
$ws.Cells.Item(21,1).Value = "ARMA_I(0,6,1)"
$ws.Cells.Item(21,7).Value = 98.12918817858983
$ws.Cells.Item(22,1).Value = "ARMA_I(0,6,2)"
$ws.Cells.Item(22,7).Value = 98.04513721315095
$ws.Cells.Item(23,1).Value = "ARMA_I(0,7,0)"
$ws.Cells.Item(23,8).Value = 97.54987876951452
$ws.Cells.Item(24,1).Value = "ARMA_I(0,7,1)"
$ws.Cells.Item(24,8).Value = 97.65553546247817
$ws.Cells.Item(25,1).Value = "ARMA_I(0,7,2)"
$ws.Cells.Item(25,8).Value = 97.58766673988303
$ws.Cells.Item(26,1).Value = "ARMA_I(1,1,0)"
$ws.Cells.Item(26,2).Value = 82.25554390024789
$ws.Cells.Item(27,1).Value = "ARMA_I(1,1,1)"
$ws.Cells.Item(27,2).Value = 86.17557709739998
$ws.Cells.Item(28,1).Value = "ARMA_I(1,10,0)"
$ws.Cells.Item(28,9).Value = 96.21438089214904
$ws.Cells.Item(29,1).Value = "ARMA_I(1,10,1)"
$ws.Cells.Item(29,9).Value = 96.14469766926345
$ws.Cells.Item(30,1).Value = "ARMA_I(1,2,0)"
$ws.Cells.Item(30,3).Value = 86.19615976744174
$ws.Cells.Item(31,1).Value = "ARMA_I(1,2,1)"
$ws.Cells.Item(31,3).Value = 86.45685250944567
$ws.Cells.Item(32,1).Value = "ARMA_I(1,3,0)"
$ws.Cells.Item(32,4).Value = 99.33277181411147
$ws.Cells.Item(33,1).Value = "ARMA_I(1,3,1)"
$ws.Cells.Item(33,4).Value = 99.30790263815177
$ws.Cells.Item(34,1).Value = "ARMA_I(1,4,0)"
$ws.Cells.Item(34,5).Value = 98.95350190457305
$ws.Cells.Item(35,1).Value = "ARMA_I(1,4,1)"
$ws.Cells.Item(35,5).Value = 98.94748754800725
$ws.Cells.Item(36,1).Value = "ARMA_I(1,5,0)"
$ws.Cells.Item(36,6).Value = 98.49284508385018
$ws.Cells.Item(37,1).Value = "ARMA_I(1,5,1)"
$ws.Cells.Item(37,6).Value = 98.55915295571171
$ws.Cells.Item(38,1).Value = "ARMA_I(1,6,0)"
$ws.Cells.Item(38,7).Value = 98.08258539137067
$ws.Cells.Item(39,1).Value = "ARMA_I(1,6,1)"
$ws.Cells.Item(39,7).Value = 98.03230335300314
$ws.Cells.Item(40,1).Value = "ARMA_I(1,7,0)"
$ws.Cells.Item(40,8).Value = 97.51807646100171
$ws.Cells.Item(41,1).Value = "ARMA_I(1,7,1)"
$ws.Cells.Item(41,8).Value = 97.54644667639212
$ws.Cells.Item(42,1).Value = "ARMA_I(2,1,0)"
$ws.Cells.Item(42,2).Value = 73.70854751749948
$ws.Cells.Item(43,1).Value = "ARMA_I(2,1,2)"
$ws.Cells.Item(43,2).Value = 81.20994996761762
$ws.Cells.Item(44,1).Value = "ARMA_I(2,10,0)"
$ws.Cells.Item(44,9).Value = 96.14283544976692
$ws.Cells.Item(45,1).Value = "ARMA_I(2,10,2)"
$ws.Cells.Item(45,9).Value = 96.05325097265769
$ws.Cells.Item(46,1).Value = "ARMA_I(2,2,0)"
$ws.Cells.Item(46,3).Value = 86.14915281209417
$ws.Cells.Item(47,1).Value = "ARMA_I(2,2,2)"
$ws.Cells.Item(47,3).Value = 86.33137278515588
$ws.Cells.Item(48,1).Value = "ARMA_I(2,3,0)"
$ws.Cells.Item(48,4).Value = 99.38644295733721
$ws.Cells.Item(49,1).Value = "ARMA_I(2,3,2)"
$ws.Cells.Item(49,4).Value = 99.3298419525742
$ws.Cells.Item(50,1).Value = "ARMA_I(2,4,0)"
$ws.Cells.Item(50,5).Value = 98.98196196311666
$ws.Cells.Item(51,1).Value = "ARMA_I(2,4,2)"
$ws.Cells.Item(51,5).Value = 98.85748858113767
$ws.Cells.Item(52,1).Value = "ARMA_I(2,5,0)"
$ws.Cells.Item(52,6).Value = 98.55592964630206
$ws.Cells.Item(53,1).Value = "ARMA_I(2,5,2)"
$ws.Cells.Item(53,6).Value = 98.40224172209965
$ws.Cells.Item(54,1).Value = "ARMA_I(2,6,0)"
$ws.Cells.Item(54,7).Value = 98.03174819738368
$ws.Cells.Item(55,1).Value = "ARMA_I(2,6,2)"
$ws.Cells.Item(55,7).Value = 97.97363221329417
$ws.Cells.Item(56,1).Value = "ARMA_I(2,7,0)"
$ws.Cells.Item(56,8).Value = 97.5259085391735
$ws.Cells.Item(57,1).Value = "ARMA_I(2,7,2)"
$ws.Cells.Item(57,8).Value = 97.57467590858005
